$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update rows 15 and 16 (two Apr-19/20-ish fixtures whose match rows were swapped back) ---
# Row 15
$ws.Range("A15").Value = 13
$ws.Range("B15").Value = 6992554
$ws.Range("C15").Value = 'Thailand Premier League'
$ws.Range("D15").Value = 'Thailand Premier League'
$ws.Range("E15").Value = 45158.33333333334
$ws.Range("F15").Value = 'Sukhothai FC'
$ws.Range("G15").Value = 'Trat FC'
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 'D'
$ws.Range("K15").Value = 1.8
$ws.Range("L15").Value = 3.6
$ws.Range("M15").Value = 4.333
$ws.Range("N15").Value = 1.833
$ws.Range("O15").Value = 3.75
$ws.Range("P15").Value = 4
$ws.Range("Q15").Value = -0.5
$ws.Range("R15").Value = 1.8
$ws.Range("S15").Value = 2
$ws.Range("T15").Value = 2.75
$ws.Range("U15").Value = 1.975
$ws.Range("V15").Value = 1.825
$ws.Range("W15").Value = -1
$ws.Range("X15").Value = 2.75
$ws.Range("Y15").Value = -1
$ws.Range("Z15").Value = -1
$ws.Range("AA15").Value = 1
$ws.Range("AB15").Value = -1
$ws.Range("AC15").Value = 0.825

# Row 16
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = 6992550
$ws.Range("C16").Value = 'Thailand Premier League'
$ws.Range("D16").Value = 'Thailand Premier League'
$ws.Range("E16").Value = 45158.33333333334
$ws.Range("F16").Value = 'Buriram United'
$ws.Range("G16").Value = 'Lamphun Warrior FC'
$ws.Range("H16").Value = 3
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 'H'
$ws.Range("K16").Value = 1.166
$ws.Range("L16").Value = 8
$ws.Range("M16").Value = 12
$ws.Range("N16").Value = 1.25
$ws.Range("O16").Value = 6.5
$ws.Range("P16").Value = 8.5
$ws.Range("Q16").Value = -1.75
$ws.Range("R16").Value = 1.875
$ws.Range("S16").Value = 1.925
$ws.Range("T16").Value = 3
$ws.Range("U16").Value = 1.975
$ws.Range("V16").Value = 1.825
$ws.Range("W16").Value = 0.25
$ws.Range("X16").Value = -1
$ws.Range("Y16").Value = -1
$ws.Range("Z16").Value = 0.875
$ws.Range("AA16").Value = -1
$ws.Range("AB16").Value = 0
$ws.Range("AC16").Value = 0

# --- Update rows 179 and 180 (match rows swapped back) ---
# Row 179
$ws.Range("A179").Value = 177
$ws.Range("B179").Value = 6992695
$ws.Range("C179").Value = 'Thailand Premier League'
$ws.Range("D179").Value = 'Thailand Premier League'
$ws.Range("E179").Value = 45385.375
$ws.Range("F179").Value = 'Muang Thong United'
$ws.Range("G179").Value = 'Uthai Thani FC'
$ws.Range("H179").Value = 5
$ws.Range("I179").Value = 2
$ws.Range("J179").Value = 'H'
$ws.Range("K179").Value = 2.1
$ws.Range("L179").Value = 3.75
$ws.Range("M179").Value = 2.7
$ws.Range("N179").Value = 1.95
$ws.Range("O179").Value = 3.8
$ws.Range("P179").Value = 2.9
$ws.Range("Q179").Value = -0.25
$ws.Range("R179").Value = 1.8
$ws.Range("S179").Value = 2
$ws.Range("T179").Value = 3
$ws.Range("U179").Value = 1.825
$ws.Range("V179").Value = 1.975
$ws.Range("W179").Value = 0.95
$ws.Range("X179").Value = -1
$ws.Range("Y179").Value = -1
$ws.Range("Z179").Value = 0.8
$ws.Range("AA179").Value = -1
$ws.Range("AB179").Value = 0.825
$ws.Range("AC179").Value = -1

# Row 180
$ws.Range("A180").Value = 178
$ws.Range("B180").Value = 8026714
$ws.Range("C180").Value = 'Thailand Premier League'
$ws.Range("D180").Value = 'Thailand Premier League'
$ws.Range("E180").Value = 45385.375
$ws.Range("F180").Value = 'BG Pathum United'
$ws.Range("G180").Value = 'Buriram United'
$ws.Range("H180").Value = 1
$ws.Range("I180").Value = 1
$ws.Range("J180").Value = 'D'
$ws.Range("K180").Value = 3
$ws.Range("L180").Value = 3.6
$ws.Range("M180").Value = 2
$ws.Range("N180").Value = 3.1
$ws.Range("O180").Value = 3.75
$ws.Range("P180").Value = 1.95
$ws.Range("Q180").Value = 0.5
$ws.Range("R180").Value = 1.825
$ws.Range("S180").Value = 1.975
$ws.Range("T180").Value = 2.75
$ws.Range("U180").Value = 1.85
$ws.Range("V180").Value = 1.95
$ws.Range("W180").Value = -1
$ws.Range("X180").Value = 2.75
$ws.Range("Y180").Value = -1
$ws.Range("Z180").Value = 0.825
$ws.Range("AA180").Value = -1
$ws.Range("AB180").Value = -1
$ws.Range("AC180").Value = 0.95

# --- Row 191: fixture now played; result + final odds filled in, id/date corrected ---
$ws.Range("A191").Value = 189
$ws.Range("B191").Value = 8075058
$ws.Range("C191").Value = 'Thailand Premier League'
$ws.Range("D191").Value = 'Thailand Premier League'
$ws.Range("E191").Value = 45400.375
$ws.Range("F191").Value = 'Bangkok United'
$ws.Range("G191").Value = 'Lamphun Warrior FC'
$ws.Range("H191").Value = 2
$ws.Range("I191").Value = 2
$ws.Range("J191").Value = 'D'
$ws.Range("K191").Value = 1.363
$ws.Range("L191").Value = 4.333
$ws.Range("M191").Value = 6.5
$ws.Range("N191").Value = 1.363
$ws.Range("O191").Value = 4.5
$ws.Range("P191").Value = 6
$ws.Range("Q191").Value = -1.5
$ws.Range("R191").Value = 1.925
$ws.Range("S191").Value = 1.775
$ws.Range("T191").Value = 3
$ws.Range("U191").Value = 1.8
$ws.Range("V191").Value = 2
$ws.Range("W191").Value = -1
$ws.Range("X191").Value = 3.5
$ws.Range("Y191").Value = -1
$ws.Range("Z191").Value = -1
$ws.Range("AA191").Value = 0.7749999999999999
$ws.Range("AB191").Value = 0.8
$ws.Range("AC191").Value = -1

# --- Row 192: update to hold the match that used to be on row 193 (BG Pathum United vs Chonburi) ---
$ws.Range("A192").Value = 190
$ws.Range("B192").Value = 6992714
$ws.Range("C192").Value = 'Thailand Premier League'
$ws.Range("D192").Value = 'Thailand Premier League'
$ws.Range("E192").Value = 45402.33333333334
$ws.Range("F192").Value = 'BG Pathum United'
$ws.Range("G192").Value = 'Chonburi'
$ws.Range("K192").Value = 1.444
$ws.Range("L192").Value = 4.333
$ws.Range("M192").Value = 6
$ws.Range("N192").Value = 1.45
$ws.Range("O192").Value = 4.5
$ws.Range("P192").Value = 5.5
$ws.Range("Q192").Value = -1.25
$ws.Range("R192").Value = 1.95
$ws.Range("S192").Value = 1.85
$ws.Range("T192").Value = 3
$ws.Range("U192").Value = 1.95
$ws.Range("V192").Value = 1.85
$ws.Range("W192").Value = 0
$ws.Range("X192").Value = 0
$ws.Range("Y192").Value = 0
$ws.Range("Z192").Value = 0
$ws.Range("AA192").Value = 0

# --- Row 193: update to hold the match that used to be on row 192 (Police Tero FC vs Uthai Thani FC) ---
$ws.Range("A193").Value = 191
$ws.Range("B193").Value = 6995900
$ws.Range("C193").Value = 'Thailand Premier League'
$ws.Range("D193").Value = 'Thailand Premier League'
$ws.Range("E193").Value = 45402.35416666666
$ws.Range("F193").Value = 'Police Tero FC'
$ws.Range("G193").Value = 'Uthai Thani FC'
$ws.Range("K193").Value = 3.3
$ws.Range("L193").Value = 3.6
$ws.Range("M193").Value = 1.95
$ws.Range("N193").Value = 3.1
$ws.Range("O193").Value = 3.6
$ws.Range("P193").Value = 2.05
$ws.Range("Q193").Value = 0.25
$ws.Range("R193").Value = 1.975
$ws.Range("S193").Value = 1.825
$ws.Range("T193").Value = 3
$ws.Range("U193").Value = 1.95
$ws.Range("V193").Value = 1.85
$ws.Range("W193").Value = 0
$ws.Range("X193").Value = 0
$ws.Range("Y193").Value = 0
$ws.Range("Z193").Value = 0
$ws.Range("AA193").Value = 0

# --- Append 3 brand-new fixture rows (194-196), copying format from row 193, then filling values ---
$ws.Range("A193:AC193").Copy($ws.Range("A194:AC196"))

# Row 194
$ws.Range("A194").Value = 192
$ws.Range("B194").Value = 6992713
$ws.Range("C194").Value = 'Thailand Premier League'
$ws.Range("D194").Value = 'Thailand Premier League'
$ws.Range("E194").Value = 45402.375
$ws.Range("F194").Value = 'Khonkaen United'
$ws.Range("G194").Value = 'Trat FC'
$ws.Range("K194").Value = 2.1
$ws.Range("L194").Value = 3.75
$ws.Range("M194").Value = 2.875
$ws.Range("N194").Value = 2.05
$ws.Range("O194").Value = 3.8
$ws.Range("P194").Value = 2.9
$ws.Range("Q194").Value = -0.25
$ws.Range("R194").Value = 1.85
$ws.Range("S194").Value = 1.95
$ws.Range("T194").Value = 3
$ws.Range("U194").Value = 2
$ws.Range("V194").Value = 1.8
$ws.Range("W194").Value = 0
$ws.Range("X194").Value = 0
$ws.Range("Y194").Value = 0
$ws.Range("Z194").Value = 0
$ws.Range("AA194").Value = 0

# Row 195
$ws.Range("A195").Value = 193
$ws.Range("B195").Value = 6992710
$ws.Range("C195").Value = 'Thailand Premier League'
$ws.Range("D195").Value = 'Thailand Premier League'
$ws.Range("E195").Value = 45402.41666666666
$ws.Range("F195").Value = 'Ratchaburi FC'
$ws.Range("G195").Value = 'Buriram United'
$ws.Range("K195").Value = 5.25
$ws.Range("L195").Value = 3.75
$ws.Range("M195").Value = 1.571
$ws.Range("N195").Value = 5
$ws.Range("O195").Value = 3.75
$ws.Range("P195").Value = 1.6
$ws.Range("Q195").Value = 1
$ws.Range("R195").Value = 1.8
$ws.Range("S195").Value = 2
$ws.Range("T195").Value = 2.75
$ws.Range("U195").Value = 1.95
$ws.Range("V195").Value = 1.85
$ws.Range("W195").Value = 0
$ws.Range("X195").Value = 0
$ws.Range("Y195").Value = 0
$ws.Range("Z195").Value = 0
$ws.Range("AA195").Value = 0

# Row 196
$ws.Range("A196").Value = 194
$ws.Range("B196").Value = 6992715
$ws.Range("C196").Value = 'Thailand Premier League'
$ws.Range("D196").Value = 'Thailand Premier League'
$ws.Range("E196").Value = 45403.375
$ws.Range("F196").Value = 'Nakhon Pathom FC'
$ws.Range("G196").Value = 'Bangkok United'
$ws.Range("K196").Value = 4.75
$ws.Range("L196").Value = 3.75
$ws.Range("M196").Value = 1.615
$ws.Range("N196").Value = 6.5
$ws.Range("O196").Value = 4
$ws.Range("P196").Value = 1.45
$ws.Range("Q196").Value = 1
$ws.Range("R196").Value = 2.025
$ws.Range("S196").Value = 1.775
$ws.Range("T196").Value = 2.75
$ws.Range("U196").Value = 1.9
$ws.Range("V196").Value = 1.9
$ws.Range("W196").Value = 0
$ws.Range("X196").Value = 0
$ws.Range("Y196").Value = 0
$ws.Range("Z196").Value = 0
$ws.Range("AA196").Value = 0
